# Insert a new weekly price record for Betarraga (Macroferia Regional de Talca)
# right before the existing row 83, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 83 (pushes old row 83..186 down to 84..187)
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new weekly record
$ws.Cells.Item(83, 1).Value = 5
$ws.Cells.Item(83, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(83, 3).Value = "Maule"
$ws.Cells.Item(83, 4).Value = 44483
$ws.Cells.Item(83, 5).Value = 7
$ws.Cells.Item(83, 6).Value = 100114014
$ws.Cells.Item(83, 7).Value = "Betarraga"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 4000
$ws.Cells.Item(83, 11).Value = 700
$ws.Cells.Item(83, 12).Value = 700
$ws.Cells.Item(83, 13).Value = 700
$ws.Cells.Item(83, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(83, 15).Value = "Región del Maule"
$ws.Cells.Item(83, 16).Value = 140
$ws.Cells.Item(83, 17).Value = 5
$ws.Cells.Item(83, 18).Value = "Hortaliza"
